$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.712.75'
$ws.Cells.Item(2, 5).Value = '  -1.99%  '

$ws.Cells.Item(3, 4).Value = '2.388.51'
$ws.Cells.Item(3, 5).Value = '  -3.25%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '549.31'
$ws.Cells.Item(5, 5).Value = '  -1.87%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '156.75'
$ws.Cells.Item(6, 5).Value = '  -3.86%  '

$ws.Cells.Item(7, 5).Value = '  -0.03%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.501'

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.155'
$ws.Cells.Item(9, 5).Value = '  +2.50%  '

$ws.Cells.Item(10, 5).Value = '  -1.60%  '

$ws.Cells.Item(11, 5).Value = '  -2.98%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.70'
$ws.Cells.Item(12, 5).Value = '  -2.41%  '

$ws.Cells.Item(13, 4).Value = '67.606.89'
$ws.Cells.Item(13, 5).Value = '  -1.75%  '

$ws.Cells.Item(14, 5).Value = '  -1.62%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '22.66'
$ws.Cells.Item(15, 5).Value = '  -4.11%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '10.23'
$ws.Cells.Item(16, 5).Value = '  -5.34%  '

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '327.29'
$ws.Cells.Item(17, 5).Value = '  -4.35%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '6.70'
$ws.Cells.Item(18, 5).Value = '  -5.50%  '

$ws.Cells.Item(19, 5).Value = '  -2.00%  '

$ws.Cells.Item(20, 5).Value = '  +0.02%  '

$ws.Cells.Item(21, 5).Value = '  -5.28%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '65.43'
$ws.Cells.Item(22, 5).Value = '  -2.63%  '

$ws.Cells.Item(23, 5).Value = '  -3.32%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '7.92'
$ws.Cells.Item(24, 5).Value = '  -3.19%  '

$ws.Cells.Item(25, 4).Value = '0.0₃0785'
$ws.Cells.Item(25, 5).Value = '  -4.04%  '

$ws.Cells.Item(26, 5).Value = '  -3.40%  '

$ws.Cells.Item(27, 5).Value = '  -0.07%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '412.88'
$ws.Cells.Item(28, 5).Value = '  -6.16%  '

$ws.Cells.Item(29, 5).Value = '  -2.90%  '

$ws.Cells.Item(30, 5).Value = '  -2.89%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '158.74'
$ws.Cells.Item(31, 5).Value = '  +1.12%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '18.96'
$ws.Cells.Item(32, 5).Value = '  -0.48%  '

$ws.Cells.Item(33, 5).Value = '  -0.06%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '17.56'
$ws.Cells.Item(34, 5).Value = '  -1.95%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.103'
$ws.Cells.Item(35, 5).Value = '  -4.61%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.289'
$ws.Cells.Item(36, 5).Value = '  -4.18%  '

$ws.Cells.Item(37, 2).Value = 'RenderToken'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '4.18'
$ws.Cells.Item(37, 5).Value = '  -6.33%  '

$ws.Cells.Item(38, 2).Value = 'Stacks'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '1.44'
$ws.Cells.Item(38, 5).Value = '  -2.61%  '

$ws.Cells.Item(39, 5).Value = '  -6.08%  '

$ws.Cells.Item(40, 2).Value = 'Filecoin'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.24'
$ws.Cells.Item(40, 5).Value = '  -3.46%  '

$ws.Cells.Item(41, 2).Value = 'Aave'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '127.03'
$ws.Cells.Item(41, 5).Value = '  -4.80%  '

$ws.Cells.Item(42, 2).Value = 'Cronos'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.0702'
$ws.Cells.Item(42, 5).Value = '  -2.42%  '

$ws.Cells.Item(43, 2).Value = 'dogwifhat'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.90'
$ws.Cells.Item(43, 5).Value = '  -8.74%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.468'
$ws.Cells.Item(44, 5).Value = '  -3.19%  '

$ws.Cells.Item(45, 5).Value = '  -2.29%  '

$ws.Cells.Item(46, 5).Value = '  -0.14%  '

$ws.Cells.Item(47, 5).Value = '  -1.20%  '

$ws.Cells.Item(48, 5).Value = '  -7.94%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '16.32'
$ws.Cells.Item(49, 5).Value = '  -3.63%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0422'
$ws.Cells.Item(50, 5).Value = '  -1.78%  '

$ws.Cells.Item(51, 4).Value = '0.0₆0198'
$ws.Cells.Item(51, 5).Value = '  -6.53%  '
